# The "Points" column (B) is being cleared out for every data row (rows 2-99),
# leaving just the header in B1. Column A (names) is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B99").ClearContents()

# Scroll the sheet so row 76 is at the top of the view, and move the active
# selection down to B99 (bottom of the data), matching where the user ended
# up after clearing the column.
$win = $excel.ActiveWindow
$win.ScrollRow = 76
$win.ScrollColumn = 1
$ws.Range("B99").Select()
